$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 197.6
$ws.Range("I33").Value = 197.6
$ws.Range("K33").Value = 197.6
$ws.Range("M33").Value = 31.40000000000001
# Row 53
$ws.Range("H53").Value = 659.6
$ws.Range("J53").Value = 504
$ws.Range("L53").Value = 504
$ws.Range("N53").Value = -1778
# Row 55
$ws.Range("H55").Value = 514.5
$ws.Range("I55").Value = 514.5
$ws.Range("K55").Value = 514.5
$ws.Range("M55").Value = -300.5
# Row 58
$ws.Range("H58").Value = 130
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 64
$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 69
$ws.Range("H69").Value = 9800
$ws.Range("J69").Value = 9800
$ws.Range("L69").Value = 29400
$ws.Range("N69").Value = -31148
# Row 72
$ws.Range("H72").Value = 9800
$ws.Range("J72").Value = 9800
$ws.Range("L72").Value = 88200
$ws.Range("N72").Value = -96936
# Row 132
$ws.Range("H132").Value = 4200
$ws.Range("I132").Value = 4200
$ws.Range("K132").Value = 12600
$ws.Range("M132").Value = -10070
# Row 138
$ws.Range("H138").Value = 4287.625
$ws.Range("I138").Value = 5498.5
$ws.Range("J138").Value = 4177.5454
$ws.Range("K138").Value = 16495.5
$ws.Range("L138").Value = 12532.6362
$ws.Range("M138").Value = -11355.5
$ws.Range("N138").Value = -22812.6362

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 2625
$ws.Range("I6").Value = 2166.6667
$ws.Range("J6").Value = 4000
$ws.Range("K6").Value = 2166.6667
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = -1993.6667
$ws.Range("N6").Value = -4346
# Row 32
$ws.Range("H32").Value = 19999.666
$ws.Range("I32").Value = 19999.666
$ws.Range("K32").Value = 19999.666
$ws.Range("M32").Value = -19712.666
# Row 45
$ws.Range("H45").Value = 7747.5
$ws.Range("I45").Value = 7747.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 7747.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -7370.5
$ws.Range("N45").ClearContents()
# Row 61
$ws.Range("H61").Value = 13193.6
$ws.Range("I61").Value = 4986
$ws.Range("K61").Value = 4986
$ws.Range("M61").Value = -4774
# Row 136
$ws.Range("H136").Value = 13193.6
$ws.Range("I136").Value = 4986
$ws.Range("K136").Value = 14958
$ws.Range("M136").Value = -12408
# Row 139
$ws.Range("H139").Value = 99997.5
$ws.Range("J139").Value = 99997.5
$ws.Range("L139").Value = 99997.5
$ws.Range("N139").Value = -110277.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1500
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1500
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1726
# Row 29
$ws.Range("H29").Value = 5628.75
$ws.Range("I29").Value = 5008
$ws.Range("J29").Value = 6249.5
$ws.Range("K29").Value = 5008
$ws.Range("L29").Value = 6249.5
$ws.Range("M29").Value = -4719
$ws.Range("N29").Value = -6827.5
# Row 36
$ws.Range("H36").Value = 6454.25
$ws.Range("I36").Value = 6454.25
$ws.Range("K36").Value = 6454.25
$ws.Range("M36").Value = -5920.25
# Row 54
$ws.Range("H54").Value = 5329.6665
$ws.Range("I54").Value = 5329.6665
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5329.6665
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4845.6665
$ws.Range("N54").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 284.75
$ws.Range("I7").Value = 284.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 284.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -171.75
$ws.Range("N7").ClearContents()
# Row 19
$ws.Range("H19").Value = 509.42856
$ws.Range("J19").Value = 724
$ws.Range("L19").Value = 724
$ws.Range("N19").Value = -1064
# Row 24
$ws.Range("H24").Value = 509.42856
$ws.Range("J24").Value = 724
$ws.Range("L24").Value = 724
$ws.Range("N24").Value = -1064
# Row 25
$ws.Range("H25").Value = 3637.3333
$ws.Range("I25").Value = 456.5
$ws.Range("J25").Value = 9999
$ws.Range("K25").Value = 456.5
$ws.Range("L25").Value = 9999
$ws.Range("M25").Value = -282.5
$ws.Range("N25").Value = -10347
# Row 31
$ws.Range("H31").Value = 9114.706
$ws.Range("I31").Value = 3090
$ws.Range("J31").Value = 11625
$ws.Range("K31").Value = 3090
$ws.Range("L31").Value = 11625
$ws.Range("M31").Value = -2795
$ws.Range("N31").Value = -12215
# Row 34
$ws.Range("H34").Value = 9114.706
$ws.Range("I34").Value = 3090
$ws.Range("J34").Value = 11625
$ws.Range("K34").Value = 3090
$ws.Range("L34").Value = 11625
$ws.Range("M34").Value = -2888
$ws.Range("N34").Value = -12029

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 4999
$ws.Range("I3").Value = 4999
$ws.Range("K3").Value = 14997
$ws.Range("M3").Value = -14885
# Row 6
$ws.Range("H6").Value = 76.5
$ws.Range("I6").Value = 76.5
$ws.Range("K6").Value = 229.5
$ws.Range("M6").Value = -116.5
# Row 12
$ws.Range("H12").Value = 50.9
$ws.Range("I12").Value = 56.5
$ws.Range("J12").Value = 47.166668
$ws.Range("K12").Value = 169.5
$ws.Range("L12").Value = 141.500004
$ws.Range("M12").Value = 3.5
$ws.Range("N12").Value = -487.500004
# Row 131
$ws.Range("H131").Value = 943.3333
$ws.Range("I131").Value = 943.3333
$ws.Range("K131").Value = 2829.9999
$ws.Range("M131").Value = 2210.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Range("H98").Value = 22089.834
$ws.Range("J98").Value = 22089.834
$ws.Range("L98").Value = 22089.834
$ws.Range("N98").Value = -28079.834
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -10224
# Row 22
$ws.Range("H22").Value = 993
$ws.Range("J22").Value = 982.5
$ws.Range("L22").Value = 982.5
$ws.Range("N22").Value = -1572.5
# Row 27
$ws.Range("H27").Value = 993
$ws.Range("J27").Value = 982.5
$ws.Range("L27").Value = 982.5
$ws.Range("N27").Value = -1196.5
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
# Row 136
$ws.Range("H136").Value = 15200.6
$ws.Range("I136").Value = 8002
$ws.Range("J136").Value = 19999.666
$ws.Range("K136").Value = 24006
$ws.Range("L136").Value = 59998.99800000001
$ws.Range("M136").Value = -21456
$ws.Range("N136").Value = -65098.99800000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 3209.2
$ws.Range("I23").Value = 273
$ws.Range("J23").Value = 5166.6665
$ws.Range("K23").Value = 273
$ws.Range("L23").Value = 5166.6665
$ws.Range("M23").Value = -44
$ws.Range("N23").Value = -5624.6665
# Row 122
$ws.Range("H122").Value = 1000
$ws.Range("J122").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("N122").Value = -7900
